# Add the six new query-result sheets (average delivery time, average order
# check, top cities, review-score histogram, top sellers, orders per year)
# to the workbook, matching the header/value layout & styling of the
# existing report sheets.

$wb = $excel.ActiveWorkbook

# Reference cell that already carries the bold/border/centered "header"
# style (style index 1) used by every existing sheet's header row.
$ws1 = $wb.Worksheets.Item(1)
$headerStyleSrc = $ws1.Range("A1")

# Reference cell that carries the plain/default style (index 0), captured
# up front (before any NumberFormat fiddling happens) so it can be used to
# strip back to "no explicit style" on text-forced cells later.
$defaultStyle = $ws1.Range("Z100").Style

function Add-ReportSheet($name) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet = $wb.Worksheets.Add($null, $lastSheet)
    $newSheet.Name = $name
    return $newSheet
}

# Writes a bold/boxed header cell, matching the style already used for
# header cells elsewhere in the workbook (copies the style only, so no new
# cellXfs entries are introduced).
function Set-HeaderCell($ws, $addr, $text) {
    $ws.Range($addr).Value = $text
    $headerStyleSrc.Copy()
    $ws.Range($addr).PasteSpecial(-4122)   # xlPasteFormats
}

# Writes a plain numeric value cell (no special style).
function Set-NumberCell($ws, $addr, $number) {
    $ws.Range($addr).Value = $number
}

# Writes a cell whose content is text that merely looks like a number
# (e.g. "-11.88"), keeping it stored as a literal string instead of being
# auto-coerced to a numeric value by the normal Value-assignment rules.
function Set-TextNumberCell($ws, $addr, $text) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($addr).Style = $defaultStyle
}

# Writes a plain text (non-numeric-looking) cell.
function Set-TextCell($ws, $addr, $text) {
    $ws.Range($addr).Value = $text
}

# ---------------------------------------------------------------------
# Sheet 5: average time of delivery
# ---------------------------------------------------------------------
$ws = Add-ReportSheet "average time of delivery"
Set-HeaderCell $ws "A1" "avg_delay_days"
Set-TextNumberCell $ws "A2" "-11.88"

# ---------------------------------------------------------------------
# Sheet 6: average order check
# ---------------------------------------------------------------------
$ws = Add-ReportSheet "average order check"
Set-HeaderCell $ws "A1" "avg_order_value"
Set-TextNumberCell $ws "A2" "137.75"

# ---------------------------------------------------------------------
# Sheet 7: top five cities by amount of cu[stomers]
# ---------------------------------------------------------------------
$ws = Add-ReportSheet "top five cities by amount of cu"
Set-HeaderCell $ws "A1" "customer_city"
Set-HeaderCell $ws "B1" "total_customers"
Set-TextCell   $ws "A2" "sao paulo"
Set-NumberCell $ws "B2" 15540
Set-TextCell   $ws "A3" "rio de janeiro"
Set-NumberCell $ws "B3" 6882
Set-TextCell   $ws "A4" "belo horizonte"
Set-NumberCell $ws "B4" 2773
Set-TextCell   $ws "A5" "brasilia"
Set-NumberCell $ws "B5" 2131
Set-TextCell   $ws "A6" "curitiba"
Set-NumberCell $ws "B6" 1521

# ---------------------------------------------------------------------
# Sheet 8: How many reviews with every sco[re]
# ---------------------------------------------------------------------
$ws = Add-ReportSheet "How many reviews with every sco"
Set-HeaderCell $ws "A1" "review_score"
Set-HeaderCell $ws "B1" "count_reviews"
Set-NumberCell $ws "A2" 1
Set-NumberCell $ws "B2" 11424
Set-NumberCell $ws "A3" 2
Set-NumberCell $ws "B3" 3151
Set-NumberCell $ws "A4" 3
Set-NumberCell $ws "B4" 8179
Set-NumberCell $ws "A5" 4
Set-NumberCell $ws "B5" 19142
Set-NumberCell $ws "A6" 5
Set-NumberCell $ws "B6" 57328

# ---------------------------------------------------------------------
# Sheet 9: Tob 5 sellers by revenue
# ---------------------------------------------------------------------
$ws = Add-ReportSheet "Tob 5 sellers by revenue"
Set-HeaderCell $ws "A1" "seller_id"
Set-HeaderCell $ws "B1" "revenue"
Set-TextCell       $ws "A2" "4869f7a5dfa277a7dca6462dcf3b52b2"
Set-TextNumberCell $ws "B2" "229472.63"
Set-TextCell       $ws "A3" "53243585a1d6dc2643021fd1853d8905"
Set-TextNumberCell $ws "B3" "222776.05"
Set-TextCell       $ws "A4" "4a3ca9315b744ce9f8e9374361493884"
Set-TextNumberCell $ws "B4" "200472.92"
Set-TextCell       $ws "A5" "fa1c13f2614d7b5c4749cbc52fecda94"
Set-TextNumberCell $ws "B5" "194042.03"
Set-TextCell       $ws "A6" "7c67e1448b00f6e969d365cea6b010ab"
Set-TextNumberCell $ws "B6" "187923.89"

# ---------------------------------------------------------------------
# Sheet 10: Number of orders by years
# ---------------------------------------------------------------------
$ws = Add-ReportSheet "Number of orders by years"
Set-HeaderCell $ws "A1" "year"
Set-HeaderCell $ws "B1" "total_orders"
Set-NumberCell $ws "A2" 2016
Set-NumberCell $ws "B2" 329
Set-NumberCell $ws "A3" 2017
Set-NumberCell $ws "B3" 45101
Set-NumberCell $ws "A4" 2018
Set-NumberCell $ws "B4" 54011
